# Insert a new price record as row 931 in the "Naranja" price sheet.
# This pushes the existing rows 931..1024 down to 932..1025 (dimension
# grows from A1:T1024 to A1:T1025), and populates the newly inserted
# row 931 with the new record's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at row 931, shifting rows 931:1024 down to 932:1025.
$ws.Range("A931:T931").Insert()

# Fill in the new record in row 931.
$ws.Cells.Item(931, 1).Value  = 5
$ws.Cells.Item(931, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(931, 3).Value  = "Maule"
$ws.Cells.Item(931, 4).Value  = 45212
$ws.Cells.Item(931, 5).Value  = 7
$ws.Cells.Item(931, 6).Value  = "Fruta"
$ws.Cells.Item(931, 7).Value  = 100102
$ws.Cells.Item(931, 8).Value  = "Cítricos"
$ws.Cells.Item(931, 9).Value  = 100102005
$ws.Cells.Item(931, 10).Value = "Naranja"
$ws.Cells.Item(931, 11).Value = "Navel Late"
$ws.Cells.Item(931, 12).Value = "Primera"
$ws.Cells.Item(931, 13).Value = 550
$ws.Cells.Item(931, 14).Value = 9000
$ws.Cells.Item(931, 15).Value = 10000
$ws.Cells.Item(931, 16).Value = 9545
$ws.Cells.Item(931, 17).Value = "`$/bandeja 15 kilos granel"
$ws.Cells.Item(931, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(931, 19).Value = 636
$ws.Cells.Item(931, 20).Value = 15
